$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data set (header + 10 data rows) after the edit described by the diff:
#  - B2 address replaced with a new one
#  - A3/A4 "Resultado" values swap (Falso positivo <-> Falso negativo)
#  - a new address row is inserted (Falso negativo / 0xd808259c...)
#  - two trailing rows (old rows 12 and 13) are removed entirely

$data = @(
    @("Resultado", "Address"),
    @("Falso positivo", "0x6d57fe045dcced8b289db59f66cd4354b6483d63"),
    @("Falso positivo", "0x1c3f580daeaac2f540c998c8ae3e4b18440f7c45"),
    @("Falso negativo", "0xd9cd7461f960e56364a294f124aac77b25e2b784"),
    @("Falso negativo", "0x684ede6645f1b71d77e0aeac519114ee8be3c410"),
    @("Falso negativo", "0x9f4562c9be26c7020909b50ccde3447f1b8c4b21"),
    @("Falso negativo", "0xd808259ca07fdf4d8fa825c4704f624352e2dc14"),
    @("Falso negativo", "0x63cfa80bbbee233a4257857dcdc9d78cbc8efe37"),
    @("Falso negativo", "0x6ddfef85ecf643628254e5af7064e05b3c6b221e"),
    @("Falso positivo", "0x44261fd8b6579ffb751569f859bc1225af0c5f2c"),
    @("Falso positivo", "0xa44643642141c2af4fb52cd320821ddfa1ad12fb")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the two now-obsolete trailing rows (old rows 12 & 13), shrinking the
# used range down to A1:B11.
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(12).Delete()
